$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 5.4
$ws.Range("G2").Value = 6.4
$ws.Range("H2").Value = 1.72
$ws.Range("I2").Value = 1.8
$ws.Range("J2").Value = 3.65
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 1.49
$ws.Range("M2").Value = 1.09
$ws.Range("N2").Value = 2.86
$ws.Range("O2").Value = 1.45
$ws.Range("P2").Value = 1.63
$ws.Range("Q2").Value = 2.28
$ws.Range("S2").Value = 4.6
$ws.Range("T2").Value = 2.16
$ws.Range("U2").Value = 1.71
$ws.Range("V2").Value = 1.97
$ws.Range("W2").Value = 1.18
$ws.Range("X2").Value = 11.5
$ws.Range("Y2").Value = 7
$ws.Range("Z2").Value = 9.800000000000001
$ws.Range("AA2").Value = 19.5
$ws.Range("AD2").Value = 11.5
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 55
$ws.Range("AG2").Value = 24
$ws.Range("AJ2").Value = 210
$ws.Range("AK2").Value = 120
$ws.Range("AL2").Value = 130
$ws.Range("AM2").Value = 260
$ws.Range("AO2").Value = 17
$ws.Range("F3").Value = 7.2
$ws.Range("G3").Value = 9
$ws.Range("H3").Value = 1.54
$ws.Range("I3").Value = 1.61
$ws.Range("J3").Value = 4.1
$ws.Range("K3").Value = 4.5
$ws.Range("L3").Value = 1.45
$ws.Range("N3").Value = 3.05
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 1.72
$ws.Range("Q3").Value = 2.14
$ws.Range("R3").Value = 1.26
$ws.Range("S3").Value = 4.1
$ws.Range("T3").Value = 2.2
$ws.Range("U3").Value = 1.68
$ws.Range("V3").Value = 2.32
$ws.Range("W3").Value = 1.11
$ws.Range("X3").Value = 980
$ws.Range("Y3").Value = 7.8
$ws.Range("Z3").Value = 10
$ws.Range("AA3").Value = 980
$ws.Range("AB3").Value = 980
$ws.Range("AC3").Value = 11.5
$ws.Range("AD3").Value = 980
$ws.Range("AE3").Value = 980
$ws.Range("AF3").Value = 75
$ws.Range("AG3").Value = 980
$ws.Range("AH3").Value = 980
$ws.Range("AI3").Value = 980
$ws.Range("AK3").Value = 180
$ws.Range("AL3").Value = 170
$ws.Range("AM3").Value = 260
$ws.Range("AO3").Value = 14
$ws.Range("G4").Value = 4.1
$ws.Range("H4").Value = 2.24
$ws.Range("I4").Value = 2.5
$ws.Range("L4").Value = 1.45
$ws.Range("N4").Value = 3.1
$ws.Range("Q4").Value = 2.1
$ws.Range("R4").Value = 1.28
$ws.Range("S4").Value = 3.85
$ws.Range("U4").Value = 1.97
$ws.Range("W4").Value = 1.34
$ws.Range("X4").Value = 14.5
$ws.Range("Y4").Value = 10.5
$ws.Range("Z4").Value = 17
$ws.Range("AA4").Value = 36
$ws.Range("AB4").Value = 14.5
$ws.Range("AC4").Value = 9.199999999999999
$ws.Range("AD4").Value = 13.5
$ws.Range("AF4").Value = 30
$ws.Range("AG4").Value = 18
$ws.Range("AH4").Value = 23
$ws.Range("AK4").Value = 60
$ws.Range("AM4").Value = 140
$ws.Range("AO4").Value = 26
$ws.Range("G5").Value = 1.66
$ws.Range("J5").Value = 3.85
$ws.Range("K5").Value = 4
$ws.Range("L5").Value = 1.5
$ws.Range("M5").Value = 1.09
$ws.Range("T5").Value = 2.12
$ws.Range("U5").Value = 1.73
$ws.Range("W5").Value = 2.5
$ws.Range("X5").Value = 12
$ws.Range("Z5").Value = 55
$ws.Range("AA5").Value = 240
$ws.Range("AB5").Value = 7
$ws.Range("AC5").Value = 8.800000000000001
$ws.Range("AE5").Value = 130
$ws.Range("AF5").Value = 8.800000000000001
$ws.Range("AG5").Value = 10.5
$ws.Range("AI5").Value = 140
$ws.Range("AJ5").Value = 16
$ws.Range("AK5").Value = 20
$ws.Range("AM5").Value = 210
$ws.Range("AN5").Value = 12.5
$ws.Range("AO5").Value = 210
$ws.Range("F6").Value = 2.28
$ws.Range("G6").Value = 2.48
$ws.Range("H6").Value = 3.75
$ws.Range("I6").Value = 4.2
$ws.Range("J6").Value = 2.92
$ws.Range("K6").Value = 3.2
$ws.Range("L6").Value = 1.63
$ws.Range("M6").Value = 1.12
$ws.Range("N6").Value = 2.38
$ws.Range("O6").Value = 1.61
$ws.Range("P6").Value = 1.45
$ws.Range("Q6").Value = 2.8
$ws.Range("S6").Value = 6
$ws.Range("T6").Value = 2.24
$ws.Range("U6").Value = 1.7
$ws.Range("V6").Value = 1.31
$ws.Range("W6").Value = 1.68
$ws.Range("X6").Value = 9.199999999999999
$ws.Range("Y6").Value = 980
$ws.Range("Z6").Value = 980
$ws.Range("AB6").Value = 8
$ws.Range("AC6").Value = 8.6
$ws.Range("AD6").Value = 980
$ws.Range("AE6").Value = 90
$ws.Range("AF6").Value = 16
$ws.Range("AG6").Value = 15
$ws.Range("AH6").Value = 980
$ws.Range("AI6").Value = 130
$ws.Range("AJ6").Value = 980
$ws.Range("AK6").Value = 980
$ws.Range("AL6").Value = 90
$ws.Range("F7").Value = 1.83
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 5.3
$ws.Range("I7").Value = 6.2
$ws.Range("J7").Value = 3.2
$ws.Range("K7").Value = 3.5
$ws.Range("M7").Value = 1.12
$ws.Range("N7").Value = 2.6
$ws.Range("O7").Value = 1.54
$ws.Range("P7").Value = 1.51
$ws.Range("Q7").Value = 2.6
$ws.Range("R7").Value = 1.18
$ws.Range("S7").Value = 5.3
$ws.Range("T7").Value = 2.28
$ws.Range("U7").Value = 1.66
$ws.Range("V7").Value = 1.19
$ws.Range("W7").Value = 2
$ws.Range("X7").Value = 8.800000000000001
$ws.Range("Y7").Value = 14.5
$ws.Range("Z7").Value = 50
$ws.Range("AA7").Value = 200
$ws.Range("AB7").Value = 6.8
$ws.Range("AC7").Value = 8.199999999999999
$ws.Range("AD7").Value = 29
$ws.Range("AF7").Value = 10.5
$ws.Range("AG7").Value = 12
$ws.Range("AH7").Value = 980
$ws.Range("AI7").Value = 170
$ws.Range("AJ7").Value = 27
$ws.Range("AK7").Value = 980
$ws.Range("AL7").Value = 70
$ws.Range("AM7").Value = 310
$ws.Range("AN7").Value = 23
